$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("summary")
$mads = $wb.Worksheets.Item("mads")
$tight = $wb.Worksheets.Item("mads_tightened")

# Update base values on the summary sheet (this cascades through formulas in
# mads / mads_tightened automatically on recalculation).
$summary.Range("B5").Value = 0.5
$summary.Range("B7").Value = 0.01
$summary.Range("B10").Value = 0.01

# Normalize the duplicate "orange" font used on a few labels in
# mads_tightened (D5, A7, A9) so it collapses onto the same font/style
# already used elsewhere (e.g. mads!A4), matching the dedup seen upstream.
$refColor = $mads.Range("A4").Font.Color
$tight.Range("D5").Font.Color = $refColor
$tight.Range("A7").Font.Color = $refColor
$tight.Range("A9").Font.Color = $refColor

# Update the selected cell (active cell) on each sheet, matching the diff.
$tight.Activate()
$tight.Range("E28").Select()

$summary.Activate()
$summary.Range("B8").Select()
